$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacità di trasmissione MW")

$ws.Range("C2").Value = 4300

$ws.Range("B3").Value = 3100
$ws.Range("D3").Value = 2900
$ws.Range("H3").Value = 300

$ws.Range("C4").Value = 2800
$ws.Range("E4").Value = 2400
$ws.Range("G4").Value = 0

$ws.Range("D5").Value = 5200

$ws.Range("E6").Value = 2400
$ws.Range("G6").Value = 1600

$ws.Range("D7").Value = 0
$ws.Range("F7").Value = 1300
$ws.Range("H7").Value = 0

$ws.Range("C8").Value = 300
$ws.Range("G8").Value = 0

$ws.Range("F9").Value = 1100
